# Daily_Scores.xlsx update
# - Revise several score values in existing rows 102-104 (2025-02-26 data).
# - Append four new rows (106-109) of scores for 2025-02-27.
# - Dimension grows from A1:Q105 to A1:Q109.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Updated values on row 102 (2025-02-26 / abs_activity)
# ---------------------------------------------------------------------------
$ws.Range("C102").Value = 7.97787821128389
$ws.Range("E102").Value = 7.35278171103401
$ws.Range("G102").Value = 9.664385131777475
$ws.Range("I102").Value = 8.914108356367102
$ws.Range("J102").Value = 8.843653859241947
$ws.Range("L102").Value = 5.646953252815712
$ws.Range("M102").Value = 9.677362108306161
$ws.Range("P102").Value = 43.58651551876864
$ws.Range("Q102").Value = 33.37484492318102

# ---------------------------------------------------------------------------
# Updated values on row 103 (2025-02-26 / rel_activity)
# ---------------------------------------------------------------------------
$ws.Range("C103").Value = 5.338073714768675
$ws.Range("M103").Value = 6.812746109883643
$ws.Range("P103").Value = 22.23681913665782

# ---------------------------------------------------------------------------
# Updated values on row 104 (2025-02-26 / abs_sleep)
# ---------------------------------------------------------------------------
$ws.Range("I104").Value = 10
$ws.Range("P104").Value = 47.96666666666667

# ---------------------------------------------------------------------------
# Helper to write a text value (date string / score-type label) without
# Excel coercing it into a date serial number.
# ---------------------------------------------------------------------------
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------------
# New row 106: 2025-02-27 / abs_activity
# ---------------------------------------------------------------------------
Set-TextValue $ws.Range("A106") "2025-02-27"
Set-TextValue $ws.Range("B106") "abs_activity"
$ws.Range("C106").Value = 9.232373959704701
$ws.Range("D106").Value = 0
$ws.Range("E106").Value = 8.562011298888885
$ws.Range("F106").Value = 0
$ws.Range("G106").Value = 10
$ws.Range("H106").Value = 10
$ws.Range("I106").Value = 10
$ws.Range("J106").Value = 10
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 5.45718708200663
$ws.Range("M106").Value = 8.171515540997287
$ws.Range("N106").Value = 0
$ws.Range("O106").Value = 0
$ws.Range("P106").Value = 45.96590079959087
$ws.Range("Q106").Value = 25.45718708200663

# ---------------------------------------------------------------------------
# New row 107: 2025-02-27 / rel_activity
# ---------------------------------------------------------------------------
Set-TextValue $ws.Range("A107") "2025-02-27"
Set-TextValue $ws.Range("B107") "rel_activity"
$ws.Range("C107").Value = 7.757413736082018
$ws.Range("D107").Value = 5
$ws.Range("E107").Value = 5.490608354359406
$ws.Range("F107").Value = 0
$ws.Range("G107").Value = 6.117991056071551
$ws.Range("H107").Value = 10
$ws.Range("I107").Value = 8.409288173643924
$ws.Range("J107").Value = 10
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 7.591512795594428
$ws.Range("N107").Value = 5
$ws.Range("O107").Value = 5
$ws.Range("P107").Value = 40.36681411575132
$ws.Range("Q107").Value = 30

# ---------------------------------------------------------------------------
# New row 108: 2025-02-27 / abs_sleep
# ---------------------------------------------------------------------------
Set-TextValue $ws.Range("A108") "2025-02-27"
Set-TextValue $ws.Range("B108") "abs_sleep"
$ws.Range("C108").Value = 7.733333333333334
$ws.Range("D108").Value = 0
$ws.Range("E108").Value = 9.933333333333334
$ws.Range("F108").Value = 0
$ws.Range("G108").Value = 9.433333333333334
$ws.Range("H108").Value = 7.133333333333333
$ws.Range("I108").Value = 8.199999999999999
$ws.Range("J108").Value = 10
$ws.Range("K108").Value = 10
$ws.Range("L108").Value = 10
$ws.Range("M108").Value = 10
$ws.Range("N108").Value = 0
$ws.Range("O108").Value = 0
$ws.Range("P108").Value = 55.3
$ws.Range("Q108").Value = 27.13333333333333

# ---------------------------------------------------------------------------
# New row 109: 2025-02-27 / rel_sleep
# ---------------------------------------------------------------------------
Set-TextValue $ws.Range("A109") "2025-02-27"
Set-TextValue $ws.Range("B109") "rel_sleep"
$ws.Range("C109").Value = 0
$ws.Range("D109").Value = 0
$ws.Range("E109").Value = 0
$ws.Range("F109").Value = 0
$ws.Range("G109").Value = 0
$ws.Range("H109").Value = 7.866311318598018
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 9.759705159705161
$ws.Range("K109").Value = 7.420877157511213
$ws.Range("L109").Value = 10
$ws.Range("M109").Value = 0
$ws.Range("N109").Value = 0
$ws.Range("O109").Value = 0
$ws.Range("P109").Value = 7.420877157511213
$ws.Range("Q109").Value = 27.62601647830318
